$d = $word.ActiveDocument
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
Write-Host "Footer tables count:" $ftr.Range.Tables.Count
$tbl = $ftr.Range.Tables(1)
$cell = $tbl.Cell(1,3)
Write-Host "Cell text:" $cell.Range.Text
